# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" for the 7deb554e-e109-4c5d-8466-7d9403c90ff4.md
# file row (row 4) across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-10-25 02:31:35"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-10-25 02:31:21"
$wsZhCn.Range("K4").Value = "2016-10-25 02:32:02"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-10-25 02:31:35"
$wsDeDe.Range("K4").Value = "2016-10-25 02:32:18"
